# Update the "precios_creatina" sheet with the newest price entry.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6's timestamp gets refreshed to a more precise value (recalculated on save).
$ws.Range("A6").Value = 45806.40671059028

# Append the new price observation as row 7.
$ws.Range("A7").Value = 45807.39295208645
$ws.Range("B7").Value = "CREATINA MONOHIDRATO EN POLVO"
$ws.Range("C7").Value = "1Kg"
$ws.Range("D7").Value = "15,41€"

# Match row 7's date formatting/style to the rest of the "fecha" column
# (copy formats only, so the value we just set is preserved).
$ws.Range("A6").Copy()
$ws.Range("A7").PasteSpecial(-4122)
